$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.314.94"
$ws.Range("E2").Value = "  +0.17%  "

$ws.Range("D3").Value = "2.279.02"
$ws.Range("E3").Value = "  -0.14%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.01"
$ws.Range("E5").Value = "  -3.85%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.90"
$ws.Range("E6").Value = "  +0.55%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.619"
$ws.Range("E7").Value = "  -0.50%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.598"
$ws.Range("E9").Value = "  -1.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.58"
$ws.Range("E10").Value = "  -2.78%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0897"
$ws.Range("E11").Value = "  -0.90%  "

$ws.Range("E12").Value = "  -1.26%  "

$ws.Range("E13").Value = "  +1.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.965"
$ws.Range("E14").Value = "  +0.25%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.02"
$ws.Range("E15").Value = "  -0.35%  "

$ws.Range("D16").Value = "2.624.89"
$ws.Range("E16").Value = "  -0.14%  "

$ws.Range("D17").Value = "2.277.21"
$ws.Range("E17").Value = "  -0.07%  "

$ws.Range("D18").Value = "42.546.84"
$ws.Range("E18").Value = "  +0.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.22"
$ws.Range("E19").Value = "  -1.68%  "

$ws.Range("E20").Value = "  -1.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.96"
$ws.Range("E21").Value = "  +1.42%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.65"
$ws.Range("E22").Value = "  -0.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.37"
$ws.Range("E23").Value = "  -7.73%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.63"
$ws.Range("E24").Value = "  -2.58%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.16"
$ws.Range("E25").Value = "  -2.39%  "

$ws.Range("E26").Value = "  -0.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.61"
$ws.Range("E27").Value = "  -2.08%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.28"
$ws.Range("E28").Value = "  -1.47%  "

$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.86"
$ws.Range("E29").Value = "  +12.98%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.03"
$ws.Range("E30").Value = "  -1.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.76"
$ws.Range("E31").Value = "  -5.79%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "164.50"
$ws.Range("E32").Value = "  +0.29%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0849"
$ws.Range("E33").Value = "  -3.09%  "

$ws.Range("E34").Value = "  -3.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.54"
$ws.Range("E35").Value = "  +0.14%  "

$ws.Range("E36").Value = "  -3.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.47"
$ws.Range("E37").Value = "  -2.64%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0345"
$ws.Range("E38").Value = "  -2.45%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.66"
$ws.Range("E39").Value = "  -0.70%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.66"
$ws.Range("E40").Value = "  -2.95%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.54"
$ws.Range("E41").Value = "  +1.78%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.49"
$ws.Range("E42").Value = "  +8.30%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.01"
$ws.Range("E43").Value = "  +0.58%  "

$ws.Range("B44").Value = "MultiversX"
$ws.Range("C44").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "68.21"
$ws.Range("E44").Value = "  -0.76%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.224"
$ws.Range("E45").Value = "  -0.50%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.84"
$ws.Range("E46").Value = "  -2.93%  "

$ws.Range("D47").Value = "1.698.82"
$ws.Range("E47").Value = "  +6.90%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "109.53"
$ws.Range("E48").Value = "  -2.78%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.74"
$ws.Range("E49").Value = "  -3.88%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.60"
$ws.Range("E50").Value = "  -3.55%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.09"
$ws.Range("E51").Value = "  -2.50%  "
